# foscarini_aumentato.xlsx edit:
#   - remove the "price" header column (it sat between "ean13" and "color")
#   - add a new "light_schema" header column right after "pic" (before "otherColors")
#
# The sheet only has data in row 1 (column headers); columns J.. shift left
# by one once "price" is deleted, and then a fresh column is opened up after
# the (now-shifted) "pic" column to host "light_schema".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Find the "price" header and delete its entire column.
$priceCell = $ws.Rows(1).Find("price")
if ($priceCell -ne $null) {
    $priceCell.EntireColumn.Delete()
}

# Find the ("now shifted") "pic" header and insert a new column right after it.
$picCell = $ws.Rows(1).Find("pic")
$newCol = $picCell.Column + 1
$ws.Cells.Item(1, $newCol).EntireColumn.Insert()
$ws.Cells.Item(1, $newCol).Value = "light_schema"
